# Daily attendance processing - 2026-01-30 09:05:47
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet
# lists the accounts that touched each attendance session, e.g.
#   "dnasr281@gmail.com, System"
# Today's processing run re-orders that list so the automated "System"
# actor is listed first, e.g.
#   "System, dnasr281@gmail.com"
# Every other "Recorded By" value (single-author rows, rows already
# ordered with System first, etc.) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$recordedByCol = 7   # column G ("Recorded By")
$lastRow = $ws.UsedRange.Rows.Count

$updated = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
        $updated++
    }
}

Write-Output "Reordered 'Recorded By' on $updated row(s)."
